$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.977.62'
$ws.Range('E2').Value = '  -2.27%  '

$ws.Range('D3').Value = '3.125.23'
$ws.Range('E3').Value = '  -0.84%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '''591.66'
$ws.Range('E5').Value = '  -3.37%  '

$ws.Range('D6').Value = '''134.99'
$ws.Range('E6').Value = '  -6.29%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = '3.116.69'
$ws.Range('E8').Value = '  -1.02%  '

$ws.Range('D9').Value = '''0.517'
$ws.Range('E9').Value = '  -1.79%  '

$ws.Range('D10').Value = '''0.145'
$ws.Range('E10').Value = '  -3.42%  '

$ws.Range('D11').Value = '''5.31'
$ws.Range('E11').Value = '  -0.94%  '

$ws.Range('D12').Value = '''0.456'
$ws.Range('E12').Value = '  -3.21%  '

$ws.Range('D13').Value = '''0.0000247'
$ws.Range('E13').Value = '  -3.38%  '

$ws.Range('D14').Value = '''33.85'
$ws.Range('E14').Value = '  -4.69%  '

$ws.Range('D15').Value = '3.630.55'
$ws.Range('E15').Value = '  -1.11%  '

$ws.Range('E16').Value = '  +1.11%  '

$ws.Range('D17').Value = '63.071.78'
$ws.Range('E17').Value = '  -2.09%  '

$ws.Range('D18').Value = '3.118.87'
$ws.Range('E18').Value = '  -1.15%  '

$ws.Range('D19').Value = '''6.67'
$ws.Range('E19').Value = '  -2.57%  '

$ws.Range('D20').Value = '''472.77'
$ws.Range('E20').Value = '  -0.89%  '

$ws.Range('D21').Value = '''14.10'
$ws.Range('E21').Value = '  -4.00%  '

$ws.Range('D22').Value = '''0.693'
$ws.Range('E22').Value = '  -4.40%  '

$ws.Range('D23').Value = '''7.60'
$ws.Range('E23').Value = '  -3.17%  '

$ws.Range('D24').Value = '''86.98'
$ws.Range('E24').Value = '  +2.72%  '

$ws.Range('D25').Value = '''12.94'
$ws.Range('E25').Value = '  -6.10%  '

$ws.Range('E26').Value = '  +0.14%  '

$ws.Range('D27').Value = '''2.70'
$ws.Range('E27').Value = '  -3.36%  '

$ws.Range('D28').Value = '''7.10'
$ws.Range('E28').Value = '  -4.03%  '

$ws.Range('D29').Value = '''7.93'
$ws.Range('E29').Value = '  -7.40%  '

$ws.Range('D30').Value = '''2.03'
$ws.Range('E30').Value = '  -3.18%  '

$ws.Range('D31').Value = '''27.04'
$ws.Range('E31').Value = '  +2.11%  '

$ws.Range('E32').Value = '  -0.02%  '

$ws.Range('E33').Value = '  -11.94%  '

$ws.Range('D34').Value = '''2.52'
$ws.Range('E34').Value = '  -5.18%  '

$ws.Range('D35').Value = '''1.08'
$ws.Range('E35').Value = '  -3.41%  '

$ws.Range('D36').Value = '''5.82'
$ws.Range('E36').Value = '  -2.18%  '

$ws.Range('D37').Value = '''51.99'
$ws.Range('E37').Value = '  -1.52%  '

$ws.Range('D38').Value = '0.0₃0709'
$ws.Range('E38').Value = '  -5.69%  '

$ws.Range('D39').Value = '''0.0387'
$ws.Range('E39').Value = '  -1.98%  '

$ws.Range('D40').Value = '''417.90'
$ws.Range('E40').Value = '  -7.78%  '

$ws.Range('E41').Value = '  -1.66%  '

$ws.Range('D42').Value = '''8.20'
$ws.Range('E42').Value = '  -1.55%  '

$ws.Range('D43').Value = '''2.69'
$ws.Range('E43').Value = '  -13.77%  '

$ws.Range('D44').Value = '2.859.67'
$ws.Range('E44').Value = '  +0.55%  '

$ws.Range('D45').Value = '''0.255'
$ws.Range('E45').Value = '  -4.90%  '

$ws.Range('E46').Value = '  -0.10%  '

$ws.Range('D47').Value = '''2.09'
$ws.Range('E47').Value = '  -8.12%  '

$ws.Range('E48').Value = '  -0.81%  '

$ws.Range('D49').Value = '''25.19'
$ws.Range('E49').Value = '  -4.66%  '

$ws.Range('D50').Value = '''2.26'
$ws.Range('E50').Value = '  -7.55%  '

$ws.Range('D51').Value = '''118.43'
$ws.Range('E51').Value = '  -2.02%  '
